$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The panel-order offsets that feed the layout geometry moved from a 0.1"
# gap to a 0.5" gap (and F3's x_start now reads off G4 instead of the old
# short-gap rule) so the rendered panels line up with the poster's actual
# seams.
$ws.Range("F3").Formula = "=G4+0.5"
$ws.Range("C4").Formula = "=D2+0.5"
$ws.Range("C5").Formula = "=D3+0.5"

# Leave the selection where the author last clicked while reviewing the
# updated panel layout.
$ws.Range("G23").Select() | Out-Null
